# Workbook: NIT-9007259935 Estado de Cuenta
# Commit: "Elimina antiguos EC y agrega nuevos y modifica Antigua BD"
# For this file, the meaningful edit is in the detail table (rows 16-20):
# the "Periodo Mora" column value moves from 2508 to 2509, and that
# column gets a center horizontal alignment it previously lacked.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$moraRange = $ws.Range("E16:E20")
$moraRange.Value = "2509"
$moraRange.HorizontalAlignment = -4108   # xlCenter
